$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Stash the original (pre-shift) pair_id text values from column A that
#    will be needed again after rows 16:17 are removed and everything below
#    shifts up two rows. Using Range.Copy (cell-to-cell) preserves the
#    original Text cell type instead of Excel re-inferring a Number from a
#    numeric-looking string (which plain .Value assignment would do).
# ---------------------------------------------------------------------------
$ws.Range("A16").Copy($ws.Range("Z1"))   # "8"
$ws.Range("A18").Copy($ws.Range("Z2"))   # "9"
$ws.Range("A20").Copy($ws.Range("Z3"))   # "10"
$ws.Range("A22").Copy($ws.Range("Z4"))   # "11"

# ---------------------------------------------------------------------------
# 2) Remove the old rows 16 and 17 (the "12G10**61993" duplicate pair) -
#    everything below shifts up by two rows, and the sheet's used range
#    shrinks from A1:X25 to A1:X23 automatically.
# ---------------------------------------------------------------------------
$ws.Rows("16:17").Delete()

# ---------------------------------------------------------------------------
# 3) Restore the renumbered pair_id text in column A for the shifted rows.
# ---------------------------------------------------------------------------
$ws.Range("Z1").Copy($ws.Range("A16"))
$ws.Range("Z1").Copy($ws.Range("A17"))
$ws.Range("Z2").Copy($ws.Range("A18"))
$ws.Range("Z2").Copy($ws.Range("A19"))
$ws.Range("Z3").Copy($ws.Range("A20"))
$ws.Range("Z3").Copy($ws.Range("A21"))
$ws.Range("Z4").Copy($ws.Range("A22"))
$ws.Range("Z4").Copy($ws.Range("A23"))

# Clean up the helper cells used for the copy round-trip.
$ws.Range("Z1:Z4").Clear()

# ---------------------------------------------------------------------------
# 4) Small corrections to existing rows (row numbers unaffected by the
#    deletion above, since they are all above row 16).
# ---------------------------------------------------------------------------
$ws.Range("G4").Value = 90
$ws.Range("O4").Value = 42856
$ws.Range("R4").Value = 136
$ws.Range("W4").Value = 17287

$ws.Range("G5").Value = 90

$ws.Range("L9").Value = "JUCA1**071955"

$ws.Range("L13").Value = "FAMU2**071985"

# ---------------------------------------------------------------------------
# 5) Correction to the (now shifted) rows 22 and 23 - id_mod typo fix.
# ---------------------------------------------------------------------------
$ws.Range("L22").Value = "JUTI1**011982"
$ws.Range("L23").Value = "JUTI1**011982"
